$wb = $excel.ActiveWorkbook

# Sheet "OFF" (Target Depth Data - Offense) - Row 3 (R / Road) values updated
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 419
$wsOff.Range("C3").Value = 293
$wsOff.Range("D3").Value = 94
$wsOff.Range("E3").Value = 41
$wsOff.Range("F3").Value = 5
$wsOff.Range("G3").Value = 6

# Sheet "DEF" (Target Depth Data - Defense) - Row 3 (R / Road) values updated
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 422
$wsDef.Range("C3").Value = 290
$wsDef.Range("D3").Value = 127
$wsDef.Range("E3").Value = 66
$wsDef.Range("F3").Value = 9
$wsDef.Range("G3").Value = 5
